$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-08 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-09 Friday", 2) | Out-Null
$d.Content.Find.Execute("322×8=2576", $true, $false, $false, $false, $false, $true, 1, $false, "483×6=2898", 2) | Out-Null
$d.Content.Find.Execute("208×9=1872", $true, $false, $false, $false, $false, $true, 1, $false, "915×3=2745", 2) | Out-Null
$d.Content.Find.Execute("333×8=2664", $true, $false, $false, $false, $false, $true, 1, $false, "567×9=5103", 2) | Out-Null
$d.Content.Find.Execute("734×6=4404", $true, $false, $false, $false, $false, $true, 1, $false, "501×9=4509", 2) | Out-Null
$d.Content.Find.Execute("443×7=3101", $true, $false, $false, $false, $false, $true, 1, $false, "837×9=7533", 2) | Out-Null
$d.Content.Find.Execute("196×3=588", $true, $false, $false, $false, $false, $true, 1, $false, "239×6=1434", 2) | Out-Null
$d.Content.Find.Execute("278×5=1390", $true, $false, $false, $false, $false, $true, 1, $false, "543×4=2172", 2) | Out-Null
$d.Content.Find.Execute("669×4=2676", $true, $false, $false, $false, $false, $true, 1, $false, "276×2=552", 2) | Out-Null
$d.Content.Find.Execute("535×8=4280", $true, $false, $false, $false, $false, $true, 1, $false, "794×9=7146", 2) | Out-Null
$d.Content.Find.Execute("827×9=7443", $true, $false, $false, $false, $false, $true, 1, $false, "635×2=1270", 2) | Out-Null
$d.Content.Find.Execute("304×3=912", $true, $false, $false, $false, $false, $true, 1, $false, "513×5=2565", 2) | Out-Null
$d.Content.Find.Execute("751×2=1502", $true, $false, $false, $false, $false, $true, 1, $false, "984×2=1968", 2) | Out-Null
$d.Content.Find.Execute("981×4=3924", $true, $false, $false, $false, $false, $true, 1, $false, "370×5=1850", 2) | Out-Null
$d.Content.Find.Execute("454×7=3178", $true, $false, $false, $false, $false, $true, 1, $false, "797×4=3188", 2) | Out-Null
$d.Content.Find.Execute("893×3=2679", $true, $false, $false, $false, $false, $true, 1, $false, "975×7=6825", 2) | Out-Null
$d.Content.Find.Execute("992×2=1984", $true, $false, $false, $false, $false, $true, 1, $false, "143×2=286", 2) | Out-Null
$d.Content.Find.Execute("777×8=6216", $true, $false, $false, $false, $false, $true, 1, $false, "691×9=6219", 2) | Out-Null
$d.Content.Find.Execute("929×9=8361", $true, $false, $false, $false, $false, $true, 1, $false, "190×5=950", 2) | Out-Null
$d.Content.Find.Execute("922×8=7376", $true, $false, $false, $false, $false, $true, 1, $false, "920×9=8280", 2) | Out-Null
$d.Content.Find.Execute("651×5=3255", $true, $false, $false, $false, $false, $true, 1, $false, "839×4=3356", 2) | Out-Null
$d.Content.Find.Execute("545×3=1635", $true, $false, $false, $false, $false, $true, 1, $false, "191×2=382", 2) | Out-Null
$d.Content.Find.Execute("865×9=7785", $true, $false, $false, $false, $false, $true, 1, $false, "139×6=834", 2) | Out-Null
$d.Content.Find.Execute("853×6=5118", $true, $false, $false, $false, $false, $true, 1, $false, "981×3=2943", 2) | Out-Null
$d.Content.Find.Execute("705×7=4935", $true, $false, $false, $false, $false, $true, 1, $false, "283×2=566", 2) | Out-Null
$d.Content.Find.Execute("746×9=6714", $true, $false, $false, $false, $false, $true, 1, $false, "696×5=3480", 2) | Out-Null
